$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MOSIP_Feature_Roadmap")

# --- Column G width change ---
$ws.Columns.Item(7).ColumnWidth = 19.7109375

# --- Append new rows 86-93 (copy formatting from row 85) ---
for ($r = 86; $r -le 93; $r++) {
  $src = $r - 1
  $ws.Range("A" + $src + ":R" + $src).Copy()
  $ws.Range("A" + $r + ":R" + $r).PasteSpecial(-4122)
  $excel.CutCopyMode = 0
}

# Fix up C and M column styles (blank cells use a different style than source row 85)
$ws.Range("A85").Copy()
$ws.Range("C86:C93").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D85").Copy()
$ws.Range("M86:M93").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row heights ---
$ws.Rows.Item(86).RowHeight = 28.5
$ws.Rows.Item(87).RowHeight = 42.75
$ws.Rows.Item(88).RowHeight = 85.5
$ws.Rows.Item(89).RowHeight = 128.25
$ws.Rows.Item(90).RowHeight = 71.25
$ws.Rows.Item(91).RowHeight = 270.75
$ws.Rows.Item(92).RowHeight = 42.75
$ws.Rows.Item(93).RowHeight = 28.5

# --- A column: running S.No. formula ---
$ws.Range("A86").Formula = "=1+A85"
$ws.Range("A87:A93").Formula = "=1+A86"

# --- B column: identification date (21-Jan-19 -> serial 43489) ---
$ws.Range("B86:B93").Value = 43489

# --- C column: blank (Reference) ---
$ws.Range("C86:C93").ClearContents()

# --- D column: Module ---
$ws.Range("D86:D93").Value = "Registration Client"

# --- E column: New / De-scoped ---
$ws.Range("E86:E93").Value = "New"

# --- F/G columns: Change Description + Logged By, row by row ---
# (interleaved to match shared-string insertion order of the target file:
#  279 Delete..., 280 Vivek/Akshaya, 281 Enter Pre-reg..., 282 Document upload..., ...)
$ws.Range("F86").Value = "Delete pre-reg packet if not consumed in client after 15 days of appointment date."
$ws.Range("G86").Value = "Vivek/Akshaya"
$ws.Range("F87").Value = "Enter Pre-registration ID: Addendum to MOS-1204:`n1. Provide ability to scan the pre-reg ID using a QR code scanner and populate the ID on the registration page."
$ws.Range("G87").Value = "Vivek/Akshaya"
$ws.Range("F88").Value = "Document upload: Addendum to MOS-1214:`n1. Document Categories + Types applicable for an individual are driven by configuration per Applicant Type + Gender + Foreigner/Local.`n2. Applicable documents are always mandatory. There is no optional document.  `n"
$ws.Range("G88").Value = "Vivek/Akshaya"
$ws.Range("F89").Value = "Preview page: Addendum to MOS-1214:`n1. Provide a timer (default 30 sec). User can proceed to the next step only after expiry of the timer.`n2. Preview page should display actual scanned images of fingerprints and irises.`n2. On navigating to Registration Preview > Edit > Modify ‘Biometric Exception’ from ‘On’ to ‘Off’ or ‘Off’ to ‘On': All biometrics previously captured (including photos) should be cleared and fresh captures will need to be made."
$ws.Range("G89").Value = "Vivek/Akshaya"
$ws.Range("F90").Value = "Acknowledgement page: Addendum to MOS-338:`nRender dummy images of left hand, right hand, thumbs, left iris and right iris. A tick or cross against each finger/Iris should indicate if the respective biometric was captured or was marked as an exception. Show fingerprint quality ranks."
$ws.Range("G90").Value = "Vivek/Akshaya"
$ws.Range("F91").Value = "UIN Update: Addendum to MOS-1299:`n1. The mandatory biometrics should be configurable. For v1 we will implement the following rules.`n1.1. UIN Update - Adult`n- Update of demographic data only: Capture at least one biometric (fingerprint or iris). More than one can be captured at the operator’s discretion.`n- Update of fingerprints: Capture all ten fingerprints minus any exceptions.`n- Update of irises: Capture both irises minus any exceptions.`n1.2. UIN Update - Child - for the first time after turning 5`n- Capture all ten fingerprints and both irises minus any exceptions.`n- Capture at least one biometric of the parent.`n1.3. UIN Update - Child - before 5 years of age`n- Only demographic data update is allowed. Capture at least one biometric of the parent.`n2. Mandatory fields to be captured and sent in the packet: UIN, List of attributes marked for update, New values of attributes, Full Name, Face photo."
$ws.Range("G91").Value = "Vivek/Akshaya"
$ws.Range("F92").Value = "Upload Packets: Addendum to MOS-559:`n1. Provide the ability for the Officer to first view packet IDs pending upload and select which ones to upload."
$ws.Range("G92").Value = "Vivek/Akshaya"
$ws.Range("F93").Value = "Device status:`nDisplay device status on the header."
$ws.Range("G93").Value = "Vivek/Akshaya"

# --- H:R columns: leave blank (clear any copied content) ---
$ws.Range("H86:R93").ClearContents()

# --- Update selection / view state ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 84
$ws.Range("C86").Select()

